$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Sema3a -> Nrp2 -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.250631
$ws.Range("H2").Value = 0.751893
$ws.Range("I2").Value = 0.2648339568266264
$ws.Range("J2").Value = 0.2648339568266264
$ws.Range("M2").Value = 25.37147633333333
$ws.Range("N2").Value = 76.114429
$ws.Range("O2").Value = 0.5780881462719274
$ws.Range("P2").Value = 0.5780881462719274
$ws.Range("Q2").Value = 6.358878484899667
$ws.Range("R2").Value = 57.229906364097
$ws.Range("S2").Value = 0.1530973711717641
$ws.Range("T2").Value = 0.1530973711717641

# Row 3 (ECs -> Sema3a -> Nrp2 -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.250631
$ws.Range("H3").Value = 0.751893
$ws.Range("I3").Value = 0.2648339568266264
$ws.Range("J3").Value = 0.2648339568266264
$ws.Range("O3").Value = 0.2328552951919536
$ws.Range("P3").Value = 0.2328552951919536
$ws.Range("Q3").Value = 2.561371542108334
$ws.Range("R3").Value = 23.052343878975
$ws.Range("S3").Value = 0.06166798919371718
$ws.Range("T3").Value = 0.06166798919371718

# Row 4 (ECs -> Sema3a -> Nrp2 -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.250631
$ws.Range("H4").Value = 0.751893
$ws.Range("I4").Value = 0.2648339568266264
$ws.Range("J4").Value = 0.2648339568266264
$ws.Range("O4").Value = 0.189056558536119
$ws.Range("P4").Value = 0.189056558536119
$ws.Range("Q4").Value = 2.079592342893333
$ws.Range("R4").Value = 18.71633108604
$ws.Range("S4").Value = 0.05006859646114509
$ws.Range("T4").Value = 0.05006859646114509

# Row 5 (FAPs -> Sema3a -> Nrp2 -> ECs)
$ws.Range("I5").Value = 0.2480790641859371
$ws.Range("J5").Value = 0.2480790641859371
$ws.Range("M5").Value = 25.37147633333333
$ws.Range("N5").Value = 76.114429
$ws.Range("O5").Value = 0.5780881462719274
$ws.Range("P5").Value = 0.5780881462719274
$ws.Range("Q5").Value = 5.956579898999556
$ws.Range("R5").Value = 53.60921909099601
$ws.Range("S5").Value = 0.1434115663441229
$ws.Range("T5").Value = 0.1434115663441229

# Row 6 (FAPs -> Sema3a -> Nrp2 -> FAPs)
$ws.Range("I6").Value = 0.2480790641859371
$ws.Range("J6").Value = 0.2480790641859371
$ws.Range("O6").Value = 0.2328552951919536
$ws.Range("P6").Value = 0.2328552951919536
$ws.Range("S6").Value = 0.05776652372195999
$ws.Range("T6").Value = 0.05776652372195999

# Row 7 (FAPs -> Sema3a -> Nrp2 -> MuSCs)
$ws.Range("I7").Value = 0.2480790641859371
$ws.Range("J7").Value = 0.2480790641859371
$ws.Range("O7").Value = 0.189056558536119
$ws.Range("P7").Value = 0.189056558536119
$ws.Range("S7").Value = 0.04690097411985424
$ws.Range("T7").Value = 0.04690097411985424

# Row 8 (MuSCs -> Sema3a -> Nrp2 -> ECs)
$ws.Range("G8").Value = 0.4609646666666666
$ws.Range("I8").Value = 0.4870869789874365
$ws.Range("J8").Value = 0.4870869789874365
$ws.Range("M8").Value = 25.37147633333333
$ws.Range("N8").Value = 76.114429
$ws.Range("O8").Value = 0.5780881462719274
$ws.Range("P8").Value = 0.5780881462719274
$ws.Range("Q8").Value = 11.69535413083622
$ws.Range("R8").Value = 105.258187177526
$ws.Range("S8").Value = 0.2815792087560404
$ws.Range("T8").Value = 0.2815792087560404

# Row 9 (MuSCs -> Sema3a -> Nrp2 -> FAPs)
$ws.Range("G9").Value = 0.4609646666666666
$ws.Range("I9").Value = 0.4870869789874365
$ws.Range("J9").Value = 0.4870869789874365
$ws.Range("O9").Value = 0.2328552951919536
$ws.Range("P9").Value = 0.2328552951919536
$ws.Range("Q9").Value = 4.710916762561111
$ws.Range("S9").Value = 0.1134207822762764
$ws.Range("T9").Value = 0.1134207822762764

# Row 10 (MuSCs -> Sema3a -> Nrp2 -> MuSCs)
$ws.Range("G10").Value = 0.4609646666666666
$ws.Range("I10").Value = 0.4870869789874365
$ws.Range("J10").Value = 0.4870869789874365
$ws.Range("O10").Value = 0.189056558536119
$ws.Range("P10").Value = 0.189056558536119
$ws.Range("S10").Value = 0.09208698795511965
$ws.Range("T10").Value = 0.09208698795511966
